$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.140.24'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.048.28'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.09%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '248.06'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').Value = '  -0.01%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '55.92'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -6.51%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.380'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -3.51%  '
$ws.Range('E10').Value = '  -2.92%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.109'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -0.27%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '16.24'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  -0.48%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.880'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  +7.56%  '
$ws.Range('D14').Value = '2.348.16'
$ws.Range('E14').Value = '  -1.33%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '5.70'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  +1.80%  '
$ws.Range('D16').Value = '2.052.19'
$ws.Range('E16').Value = '  -1.37%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '18.39'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +13.71%  '
$ws.Range('D18').Value = '37.141.32'
$ws.Range('E18').Value = '  -0.49%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '74.49'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = '0.0₃0892'
$ws.Range('E20').Value = '  -4.03%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '5.39'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -2.12%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '236.51'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -1.32%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +2.23%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '9.53'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +1.47%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '169.53'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  -5.74%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '20.03'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -2.01%  '
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  -0.96%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '4.84'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +1.07%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '0.0618'
$cell.Style = "Normal"
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '4.48'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -0.96%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.0888'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -2.86%  '
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -3.24%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '5.27'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +14.41%  '
$ws.Range('E40').Value = '  +9.17%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.0987'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -17.06%  '
$ws.Range('E42').Value = '  -2.42%  '
$ws.Range('E43').Value = '  -2.06%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '17.25'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -3.94%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '95.44'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('E46').Value = '  -3.48%  '
$ws.Range('D47').Value = '1.266.56'
$ws.Range('E47').Value = '  -3.41%  '
$ws.Range('E48').Value = '  -3.15%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '6.77'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').Value = '2.229.29'
$ws.Range('E50').Value = '  -1.63%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '43.78'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -1.79%  '
